$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F43").Value = 16
$ws.Range("G43").Value = 3086.24
$ws.Range("F48").Value = 19
$ws.Range("G48").Value = 741.76
$ws.Range("B63").Value = 38459.44
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 0
$ws.Range("F104").Value = 26
$ws.Range("G104").Value = 9624.68
$ws.Range("F113").Value = 4
$ws.Range("G113").Value = 539.08
$ws.Range("F119").Value = 7
$ws.Range("G119").Value = 792.8200000000001
$ws.Range("F121").Value = 120
$ws.Range("G121").Value = 2335.2
$ws.Range("B122").Value = 242331.18
$ws.Range("F148").Value = 37
$ws.Range("G148").Value = 1830.76
$ws.Range("F155").Value = 112
$ws.Range("G155").Value = 2204.16
$ws.Range("F157").Value = 139
$ws.Range("G157").Value = 6188.28
$ws.Range("F158").Value = 26
$ws.Range("G158").Value = 832.52
$ws.Range("B160").Value = 21842.01
$ws.Range("B165").Value = 57756
$ws.Range("B166").Value = 53925
$ws.Range("F167").Value = 3
$ws.Range("G167").Value = 137.7
$ws.Range("B168").Value = 562.75
$ws.Range("F198").Value = 19
$ws.Range("G198").Value = 2293.11
$ws.Range("B206").Value = 16535.36
$ws.Range("F210").Value = 130
$ws.Range("G210").Value = 4033.9
$ws.Range("B218").Value = 13717.7
$ws.Range("F223").Value = 66
$ws.Range("G223").Value = 3532.98
$ws.Range("B225").Value = 5126.55
$ws.Range("F235").Value = 24
$ws.Range("G235").Value = 2035.92
$ws.Range("B238").Value = 8673.75
$ws.Range("F253").Value = 1
$ws.Range("G253").Value = 476.75
$ws.Range("B259").Value = 2164.7
$ws.Range("F303").Value = 84
$ws.Range("G303").Value = 18162.48
$ws.Range("F309").Value = 17
$ws.Range("G309").Value = 1409.98
$ws.Range("F317").Value = 92
$ws.Range("G317").Value = 12610.44
$ws.Range("F323").Value = 142
$ws.Range("G323").Value = 7287.44
$ws.Range("F332").Value = 196
$ws.Range("G332").Value = 10217.48
$ws.Range("F333").Value = 148
$ws.Range("G333").Value = 16570.08
$ws.Range("F334").Value = 156
$ws.Range("G334").Value = 20125.56
$ws.Range("F341").Value = 16
$ws.Range("G341").Value = 1420.8
$ws.Range("F350").Value = 36
$ws.Range("G350").Value = 4270.68
$ws.Range("F351").Value = 69
$ws.Range("G351").Value = 4079.97
$ws.Range("F361").Value = 16
$ws.Range("G361").Value = 2975.84
$ws.Range("B370").Value = 345011.02
$ws.Range("F423").Value = 32
$ws.Range("G423").Value = 3091.2
$ws.Range("B428").Value = 46400.2
$ws.Range("F433").Value = 52
$ws.Range("G433").Value = 5109
$ws.Range("F437").Value = 50
$ws.Range("G437").Value = 1871
$ws.Range("F438").Value = 58
$ws.Range("G438").Value = 10813.52
$ws.Range("F444").Value = 62
$ws.Range("G444").Value = 3391.4
$ws.Range("B445").Value = 42096.49
$ws.Range("F460").Value = 3
$ws.Range("G460").Value = 11355.9
$ws.Range("B466").Value = 91528.91
$ws.Range("F470").Value = 606
$ws.Range("G470").Value = 7762.86
$ws.Range("F472").Value = 205
$ws.Range("G472").Value = 4044.65
$ws.Range("F475").Value = 375
$ws.Range("G475").Value = 7398.75
$ws.Range("F476").Value = 427
$ws.Range("G476").Value = 2809.66
$ws.Range("F481").Value = 954
$ws.Range("G481").Value = 6191.46
$ws.Range("F482").Value = 435
$ws.Range("G482").Value = 5720.25
$ws.Range("F485").Value = 614
$ws.Range("G485").Value = 9044.219999999999
$ws.Range("B486").Value = 104749.93
$ws.Range("F520").Value = 6
$ws.Range("G520").Value = 3976.62
$ws.Range("B524").Value = 20651.03
$ws.Range("F535").Value = 137
$ws.Range("G535").Value = 3614.06
$ws.Range("B537").Value = 40224.33
$ws.Range("F566").Value = 58
$ws.Range("G566").Value = 5713
$ws.Range("F567").Value = 106
$ws.Range("G567").Value = 10441
$ws.Range("F569").Value = 42
$ws.Range("G569").Value = 4137
$ws.Range("B575").Value = 67752
$ws.Range("F584").Value = 18
$ws.Range("G584").Value = 60.84
$ws.Range("B587").Value = 6881.2
$ws.Range("F602").Value = 67
$ws.Range("G602").Value = 8229.610000000001
$ws.Range("B604").Value = 37486.75
$ws.Range("F610").Value = 74
$ws.Range("G610").Value = 2012.8
$ws.Range("B613").Value = 62699.19
$ws.Range("F639").Value = 170
$ws.Range("G639").Value = 7340.6
$ws.Range("B641").Value = 29160.56
$ws.Range("F713").Value = 42
$ws.Range("G713").Value = 1740.9
$ws.Range("F714").Value = 11
$ws.Range("G714").Value = 591.03
$ws.Range("F715").Value = 21
$ws.Range("G715").Value = 2412.06
$ws.Range("B721").Value = 480989.07
$ws.Range("F732").Value = 38
$ws.Range("G732").Value = 3910.58
$ws.Range("B739").Value = 19605.36
$ws.Range("F742").Value = 30
$ws.Range("G742").Value = 1122
$ws.Range("F744").Value = 81
$ws.Range("G744").Value = 3029.4
$ws.Range("F746").Value = 86
$ws.Range("G746").Value = 3216.4
$ws.Range("B747").Value = 7630.18
$ws.Range("F789").Value = 4
$ws.Range("G789").Value = 119.68
$ws.Range("B790").Value = 2361.67
$ws.Range("F792").Value = 838
$ws.Range("G792").Value = 136686.18
$ws.Range("F796").Value = 75
$ws.Range("G796").Value = 5062.5
$ws.Range("B797").Value = 163463.92
$ws.Range("B803").Value = 3044419.19
$ws.Range("B804").Value = 3044419.19
